# A new weekly price record was inserted at row 282 of the data table.
# This shifts all existing records from row 282 onward down by one row
# (the former row 397 becomes row 398), and populates the newly opened
# row 282 with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 282, pushing rows 282:397 down to 283:398.
$ws.Rows("282:282").Insert()

# Populate the new row 282 with the new "Zapallo italiano" price record.
$ws.Cells.Item(282, 1).Value  = 5
$ws.Cells.Item(282, 2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(282, 3).Value  = "Maule"
$ws.Cells.Item(282, 4).Value  = 44784
$ws.Cells.Item(282, 5).Value  = 7
$ws.Cells.Item(282, 6).Value  = 100112032
$ws.Cells.Item(282, 7).Value  = "Zapallo italiano"
$ws.Cells.Item(282, 8).Value  = "Sin especificar"
$ws.Cells.Item(282, 9).Value  = "Primera"
$ws.Cells.Item(282, 10).Value = 300
$ws.Cells.Item(282, 11).Value = 20000
$ws.Cells.Item(282, 12).Value = 20000
$ws.Cells.Item(282, 13).Value = 20000
$ws.Cells.Item(282, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(282, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(282, 16).Value = 400
$ws.Cells.Item(282, 17).Value = 50
$ws.Cells.Item(282, 18).Value = "Hortaliza"
